$wb = $excel.ActiveWorkbook

# --- Update "Metadata" sheet: refresh the Last Updated timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 09:35 AM"

# --- Update "Stock List" sheet: the data window rolled forward by one row ---
# Row 2 (CAPTRU-RE1) dropped off the top, every remaining row shifted up by
# one, and a new row was appended at the bottom for TRAVELFOOD.
$ws = $wb.Worksheets.Item("Stock List")
$ws.Rows.Item(2).Delete()

$ws.Range("A76").Value = "📋"
$ws.Range("B76").Value = "TRAVELFOOD"
$ws.Range("C76").Value = "TRAVELFOOD"
$ws.Range("D76").Value = 1316.3
$ws.Range("E76").Value = 0.1141
$ws.Range("F76").Value = "N/A"
$ws.Range("G76").Value = "N/A"
$ws.Range("H76").Value = 17332.9705
